# Update the "Förändrad" (Changed) date column (C) from 45205 to 45206
# for all data rows (rows 2 through 289) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 289

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
